$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.236.52"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "'1.860.00"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'0.7001"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.08168"
$ws.Range("E8").Value = "  +9.22%  "
$ws.Range("D9").Value = "'0.3033"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'23.20"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'0.08157"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'1.862.43"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'5.157"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'0.7101"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "'89.07"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'29.255.16"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "'5.773"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.33"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007819"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "'236.13"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'2.113.01"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'7.453"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("D25").Value = "'161.57"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'8.950"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").Value = "'0.1436"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "'1.961"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'1.430"
$ws.Range("E30").Value = "  +3.01%  "
$ws.Range("D31").Value = "'1.479"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'4.382"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").Value = "'4.059"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").Value = "'0.05193"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "'1.167"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "'0.7067"
$ws.Range("D37").Value = "'0.9989"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").Value = "'2.674"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "'0.01839"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("D40").Value = "'2.726"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9230"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.137.21"
$ws.Range("E42").Value = "  +4.88%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4269"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.908"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'70.41"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'102.23"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'1.768"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("D49").Value = "'2.009.06"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'9.155"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "'6.942"
$ws.Range("E51").Value = "  -1.37%  "
